$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title heading
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our review of Acorn Pixie, a 5-reel slot game with 30 to 50 paylines by Bally Technologies. Try it free and enjoy stunning graphics and bonuses.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the document
#    and 3) replace the italic paragraph's text with the new DALLE prompt
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($count - 1)
$boldTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)
$italicRange = $italicPara.Range
$italicTextRange = $d.Range($italicRange.Start, $italicRange.End - 1)
$italicTextRange.Text = "Prompt for DALLE: Create a feature image for Acorn Pixie that showcases a happy Maya warrior in cartoon style with glasses. The image should incorporate elements from the game such as fairies, elves, flowers, and mushrooms, with a mystical forest backdrop. The Maya warrior can be holding an acorn or surrounded by them, with the game's logo appearing somewhere in the image. Make the image bright, colorful, and enticing to capture the attention of potential players."
